$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Set the data rows (A2:H11) with the refreshed listing snapshot.
# Row 2 is a newly observed listing; rows 3-10 are the previous rows 2-9
# shifted down by one with refreshed timestamps; row 11 is a newly
# appended listing.
$ws.Cells.Item(2, 1).Value = '2026-01-17 18:26:13'
$ws.Cells.Item(2, 2).Value = '製造業向け設備要件定義書の自動生成AIシステムの開発・DB設計支援エンジニア(AI/バックエンド)'
$ws.Cells.Item(2, 3).Value = 'システム開発'
$ws.Cells.Item(2, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(2, 5).Value = '期限情報なし'
$ws.Cells.Item(2, 6).Value = 'https://www.lancers.jp/work/detail/5473648'
$ws.Cells.Item(2, 7).Value = 390
$ws.Cells.Item(2, 8).Value = '🔥AI,Ai ◆開発'

$ws.Cells.Item(3, 1).Value = '2026-01-17 18:26:13'
$ws.Cells.Item(3, 2).Value = '【急募】airtableで社内業務管理システムを共に構築してくれる方'
$ws.Cells.Item(3, 3).Value = 'システム開発'
$ws.Cells.Item(3, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(3, 5).Value = '期限情報なし'
$ws.Cells.Item(3, 6).Value = 'https://www.lancers.jp/work/detail/5473383'
$ws.Cells.Item(3, 7).Value = 353
$ws.Cells.Item(3, 8).Value = '🔥AI,Ai ◇管理'

$ws.Cells.Item(4, 1).Value = '2026-01-17 18:26:13'
$ws.Cells.Item(4, 2).Value = '【シンプル版】生成AIデジタル・コミュニティ制作の依頼'
$ws.Cells.Item(4, 3).Value = 'システム開発'
$ws.Cells.Item(4, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(4, 5).Value = '期限情報なし'
$ws.Cells.Item(4, 6).Value = 'https://www.lancers.jp/work/detail/5469128'
$ws.Cells.Item(4, 7).Value = 310
$ws.Cells.Item(4, 8).Value = '🔥AI,Ai'

$ws.Cells.Item(5, 1).Value = '2026-01-17 18:26:13'
$ws.Cells.Item(5, 2).Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Cells.Item(5, 3).Value = 'システム開発'
$ws.Cells.Item(5, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(5, 5).Value = '期限情報なし'
$ws.Cells.Item(5, 6).Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Cells.Item(5, 7).Value = 243
$ws.Cells.Item(5, 8).Value = '🔥API ◆ツール'

$ws.Cells.Item(6, 1).Value = '2026-01-17 18:26:13'
$ws.Cells.Item(6, 2).Value = '※急募:Next.jsによる業務アプリの開発(+Flutter)'
$ws.Cells.Item(6, 3).Value = 'システム開発'
$ws.Cells.Item(6, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(6, 5).Value = '期限情報なし'
$ws.Cells.Item(6, 6).Value = 'https://www.lancers.jp/work/detail/5473147'
$ws.Cells.Item(6, 7).Value = 225
$ws.Cells.Item(6, 8).Value = '🔥Next.js ◆開発 ◇アプリ'

$ws.Cells.Item(7, 1).Value = '2026-01-17 18:26:13'
$ws.Cells.Item(7, 2).Value = '※急募:Flutterによる業務アプリの開発(+next.js)'
$ws.Cells.Item(7, 3).Value = 'システム開発'
$ws.Cells.Item(7, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(7, 5).Value = '期限情報なし'
$ws.Cells.Item(7, 6).Value = 'https://www.lancers.jp/work/detail/5473146'
$ws.Cells.Item(7, 7).Value = 218
$ws.Cells.Item(7, 8).Value = '🔥Next.js ◆開発 ◇アプリ'

$ws.Cells.Item(8, 1).Value = '2026-01-17 18:26:13'
$ws.Cells.Item(8, 2).Value = '【急募】Accessでの受発注管理・請求書発行システム開発'
$ws.Cells.Item(8, 3).Value = 'システム開発'
$ws.Cells.Item(8, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(8, 5).Value = '期限情報なし'
$ws.Cells.Item(8, 6).Value = 'https://www.lancers.jp/work/detail/5473234'
$ws.Cells.Item(8, 7).Value = 148
$ws.Cells.Item(8, 8).Value = '◆開発,システム開発 ◇管理'

$ws.Cells.Item(9, 1).Value = '2026-01-17 18:26:13'
$ws.Cells.Item(9, 2).Value = '【バイナリ解析 / 逆コンパイル】EPCデータ解析ツール開発'
$ws.Cells.Item(9, 3).Value = 'システム開発'
$ws.Cells.Item(9, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(9, 5).Value = '期限情報なし'
$ws.Cells.Item(9, 6).Value = 'https://www.lancers.jp/work/detail/5473394'
$ws.Cells.Item(9, 7).Value = 135
$ws.Cells.Item(9, 8).Value = '◆ツール,開発'

$ws.Cells.Item(10, 1).Value = '2026-01-17 18:26:13'
$ws.Cells.Item(10, 2).Value = '【バイナリ解析 / 逆コンパイル】EPCデータ解析ツール開発'
$ws.Cells.Item(10, 3).Value = 'システム開発'
$ws.Cells.Item(10, 4).Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Cells.Item(10, 5).Value = '期限情報なし'
$ws.Cells.Item(10, 6).Value = 'https://www.lancers.jp/work/detail/5473181'
$ws.Cells.Item(10, 7).Value = 135
$ws.Cells.Item(10, 8).Value = '◆ツール,開発'

$ws.Cells.Item(11, 1).Value = '2026-01-17 18:26:13'
$ws.Cells.Item(11, 2).Value = '製造業DXプロダクト開発のプロダクトマネージャー募集'
$ws.Cells.Item(11, 3).Value = 'システム開発'
$ws.Cells.Item(11, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(11, 5).Value = '期限情報なし'
$ws.Cells.Item(11, 6).Value = 'https://www.lancers.jp/work/detail/5468432'
$ws.Cells.Item(11, 7).Value = 75
$ws.Cells.Item(11, 8).Value = '◆開発'

# Rebuild the hyperlinks on column F (URL) so each link points at the
# right row after the reshuffle, reusing the workbook's Hyperlink style.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), 'https://www.lancers.jp/work/detail/5473648') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), 'https://www.lancers.jp/work/detail/5473383') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), 'https://www.lancers.jp/work/detail/5469128') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), 'https://www.lancers.jp/work/detail/5217096') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), 'https://www.lancers.jp/work/detail/5473147') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), 'https://www.lancers.jp/work/detail/5473146') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), 'https://www.lancers.jp/work/detail/5473234') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), 'https://www.lancers.jp/work/detail/5473394') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), 'https://www.lancers.jp/work/detail/5473181') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), 'https://www.lancers.jp/work/detail/5468432') | Out-Null
$ws.Range("F2:F11").Style = "Hyperlink"

# Widen column B (title) slightly to fit the new, longer title text.
$ws.Columns.Item(2).ColumnWidth = 50.166666666666664

